# TradingModel - 2021/11/15 watching updated
#
# The workbook is the "OpenPositionWatching" sheet that a Python job
# refreshes daily: the oldest row drops off, every remaining row's
# numbers move, and a brand new row/column pair of bars appears with
# this run's figures. We just push the freshly computed numbers into
# the sheet and refresh the view the way the author's Excel session
# left it (cursor on F11, and the Previous_* columns widened to fit
# their header text).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 2 ----
$ws.Range("B2").Value = 2436
$ws.Range("C2").Value = 60
$ws.Range("D2").Value = 204.8
$ws.Range("E2").Value = 6020
$ws.Range("F2").Value = 88.9
$ws.Range("G2").Value = 66.8
$ws.Range("H2").Value = 88.9
$ws.Range("I2").Value = 66.8

# ---- Row 3 ----
$ws.Range("B3").Value = 3035
$ws.Range("C3").Value = 32
$ws.Range("D3").Value = 185
$ws.Range("E3").Value = 5920
$ws.Range("F3").Value = 121.5
$ws.Range("G3").Value = 86.7
$ws.Range("H3").Value = 176
$ws.Range("I3").Value = 163

# ---- Row 4 ----
$ws.Range("B4").Value = 3122
$ws.Range("C4").Value = 90
$ws.Range("D4").Value = 131
$ws.Range("E4").Value = 5865
$ws.Range("F4").Value = 61.7
$ws.Range("G4").Value = 42.45
$ws.Range("H4").Value = 61.7
$ws.Range("I4").Value = 42.45

# ---- Row 5 ----
$ws.Range("B5").Value = 3141
$ws.Range("C5").Value = 27
$ws.Range("D5").Value = 218.5
$ws.Range("E5").Value = 5899.5
$ws.Range("F5").Value = 185
$ws.Range("G5").Value = 127
$ws.Range("H5").Value = 237.5
$ws.Range("I5").Value = 213.5

# ---- Row 6 ----
$ws.Range("B6").Value = 3588
$ws.Range("C6").Value = 35
$ws.Range("D6").Value = 155
$ws.Range("E6").Value = 5425
$ws.Range("F6").Value = 142.5
$ws.Range("G6").Value = 111
$ws.Range("H6").Value = 149
$ws.Range("I6").Value = 111

# ---- Row 7 ----
$ws.Range("B7").Value = 6104
$ws.Range("C7").Value = 36
$ws.Range("D7").Value = 170
$ws.Range("E7").Value = 6120
$ws.Range("F7").Value = 155
$ws.Range("G7").Value = 113.5
$ws.Range("H7").Value = 160
$ws.Range("I7").Value = 113.5

# ---- Row 8 ----
$ws.Range("B8").Value = 6138
$ws.Range("C8").Value = 30
$ws.Range("D8").Value = 203
$ws.Range("E8").Value = 6090
$ws.Range("F8").Value = 193
$ws.Range("G8").Value = 144
$ws.Range("H8").Value = 193
$ws.Range("I8").Value = 176

# ---- Row 9 ----
$ws.Range("B9").Value = 6271
$ws.Range("C9").Value = 20
$ws.Range("D9").Value = 302.5
$ws.Range("E9").Value = 6050
$ws.Range("F9").Value = 279.5
$ws.Range("G9").Value = 215.5
$ws.Range("H9").Value = 288
$ws.Range("I9").Value = 275

# ---- Row 10 (only the Previous_Platform_High/Low pair changed here) ----
$ws.Range("F10").Value = 222.5
$ws.Range("G10").Value = 180

# Widen the Previous_Platform_High/Low/N_High/N_Low columns so their
# (longer) header text fits, same as Excel's own "AutoFit Column Width"
# would leave behind.
$ws.Columns("F").ColumnWidth = 26.428571428571427
$ws.Columns("G").ColumnWidth = 25.857142857142858
$ws.Columns("H").ColumnWidth = 18.857142857142858
$ws.Columns("I").ColumnWidth = 18.285714285714285

# Leave the selection where the author's session ended up.
$ws.Range("F11").Select()
